$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.450.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.34%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.314.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.44%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.64%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.12%  "

# Row 7
$ws.Range("E7").Value = "  +0.14%  "

# Row 8
$ws.Range("E8").Value = "  -3.18%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.310.87"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.34%  "

# Row 10
$ws.Range("E10").Value = "  -6.07%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.573"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.33%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.59%  "

# Row 13
$ws.Range("E13").Value = "  -5.40%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "666.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.74%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.855.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.88%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.67%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.603.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.93%  "

# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.332.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.55%  "

# Row 19
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.118"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.27%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.59%  "

# Row 21
$ws.Range("E21").Value = "  -2.79%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.885"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.45%  "

# Row 23
$ws.Range("E23").Value = "  +4.69%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.39%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.64%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.69%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.25%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.63%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.85%  "

# Row 30
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.59%  "

# Row 31
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.16%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "581.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.65%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.45%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.103"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.42%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.04%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.703.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.38%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.20%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -16.82%  "

# Row 39
$ws.Range("E39").Value = "  -0.46%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "32.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.82%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.34%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.54%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.330"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.59%  "

# Row 44
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.94%  "

# Row 45
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₃0656"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.66%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0404"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.01%  "

# Row 47
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.127"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.67%  "

# Row 48
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.02%  "

# Row 49
$ws.Range("E49").Value = "  +0.14%  "

# Row 50
$ws.Range("E50").Value = "  -4.51%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "126.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.81%  "
